# Matriz RACI - avances para revisión 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12's activity was re-scoped from a generic "specialist users" pilot
# to a pilot specifically with deaf users.
$ws.Range("A12").Value = "Validación piloto con usuarios sordos"

# Leave the selection where the author was working when they saved.
$ws.Range("C20").Select()
